# Apply the cryptos-list price/volume refresh (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.472.54'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.84%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.146.79'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.32%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '606.17'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.98%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.84'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.51%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.143.14'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.32%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.524'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.22%  '
$ws.Range("E10").Value = '  +0.77%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.42'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.37%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.467'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.34%  '
$ws.Range("E13").Value = '  +3.40%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.42'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.50%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.666.46'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.37%  '
$ws.Range("E16").Value = '  +2.83%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.455.31'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.77%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.148.05'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.33%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.86'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.38%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '480.23'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.01%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.64'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.03%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.715'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.73%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.68'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.77%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '85.04'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.00%  '
$ws.Range("E25").Value = '  -0.63%  '
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("E27").Value = '  -1.27%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.42'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.23%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.21'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +7.35%  '
$ws.Range("E30").Value = '  +2.67%  '
$ws.Range("E31").Value = '  -4.87%  '
$ws.Range("B32").Value = 'FirstDigitalUSD'
$ws.Range("C32").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.00'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.03%  '
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.96'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.23%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.65'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.85%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.10'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.31%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.98'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.69%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0₃0774'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.84%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '52.39'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.47%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.02'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.87%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '447.35'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.59%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0396'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.92%  '
$ws.Range("E42").Value = '  +1.24%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.23'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.49%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.862.82'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.94%  '
$ws.Range("E45").Value = '  -0.87%  '
$ws.Range("E46").Value = '  -0.71%  '
$ws.Range("E47").Value = '  +2.13%  '
$ws.Range("E48").Value = '  -0.02%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '26.16'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.42%  '
$ws.Range("E50").Value = '  -0.27%  '
$ws.Range("E51").Value = '  +1.09%  '
